$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1647.6
$ws.Range("J17").Value = 2188.4614
$ws.Range("L17").Value = 6565.3842
$ws.Range("N17").Value = -6901.3842
$ws.Range("H28").Value = 525
$ws.Range("J28").Value = 303
$ws.Range("L28").Value = 303
$ws.Range("N28").Value = -1273
$ws.Range("H100").Value = 2052.7222
$ws.Range("I100").Value = 2012
$ws.Range("K100").Value = 2012
$ws.Range("M100").Value = -1471
$ws.Range("H131").Value = 30714.285
$ws.Range("I131").Value = 2500
$ws.Range("K131").Value = 7500
$ws.Range("M131").Value = -2460
$ws.Range("H137").Value = 20413526
$ws.Range("J137").Value = 7117.0303
$ws.Range("L137").Value = 21351.0909
$ws.Range("N137").Value = -26451.0909

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 576.4231
$ws.Range("I2").Value = 457.79166
$ws.Range("K2").Value = 457.79166
$ws.Range("M2").Value = -344.79166
$ws.Range("H32").Value = 130608.17
$ws.Range("I32").Value = 188951.39
$ws.Range("K32").Value = 188951.39
$ws.Range("M32").Value = -188664.39
$ws.Range("H61").Value = 1496111.5
$ws.Range("I61").Value = 3573.8113
$ws.Range("K61").Value = 3573.8113
$ws.Range("M61").Value = -3361.8113
$ws.Range("H102").Value = 2086
$ws.Range("I102").Value = 2216.1538
$ws.Range("J102").Value = 1240
$ws.Range("K102").Value = 2216.1538
$ws.Range("L102").Value = 1240
$ws.Range("M102").Value = -594.1538
$ws.Range("N102").Value = -4484
$ws.Range("H116").Value = 576.4231
$ws.Range("I116").Value = 457.79166
$ws.Range("K116").Value = 457.79166
$ws.Range("M116").Value = 1836.20834
$ws.Range("H136").Value = 1496111.5
$ws.Range("I136").Value = 3573.8113
$ws.Range("K136").Value = 10721.4339
$ws.Range("M136").Value = -8171.4339

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 576.4231
$ws.Range("I3").Value = 457.79166
$ws.Range("K3").Value = 457.79166
$ws.Range("M3").Value = -343.79166
$ws.Range("H20").Value = 35545.066
$ws.Range("I20").Value = 47237
$ws.Range("K20").Value = 47237
$ws.Range("M20").Value = -46990
$ws.Range("H86").Value = 1860.9
$ws.Range("J86").Value = 2019.4
$ws.Range("L86").Value = 2019.4
$ws.Range("N86").Value = -4265.4
$ws.Range("H89").Value = 1860.9
$ws.Range("J89").Value = 2019.4
$ws.Range("L89").Value = 10097
$ws.Range("N89").Value = -21329
$ws.Range("H94").Value = 1678.7587
$ws.Range("I94").Value = 1759.5454
$ws.Range("J94").Value = 1424.8572
$ws.Range("K94").Value = 1759.5454
$ws.Range("L94").Value = 1424.8572
$ws.Range("M94").Value = -1308.5454
$ws.Range("N94").Value = -2326.8572
$ws.Range("H99").Value = 7605.2354
$ws.Range("I99").Value = 10717.182
$ws.Range("J99").Value = 1900
$ws.Range("K99").Value = 10717.182
$ws.Range("L99").Value = 1900
$ws.Range("M99").Value = -9219.182000000001
$ws.Range("N99").Value = -4896
$ws.Range("H105").Value = 9322.223
$ws.Range("I105").Value = 8384
$ws.Range("K105").Value = 8384
$ws.Range("M105").Value = -6637
$ws.Range("H134").Value = 4769987
$ws.Range("I134").Value = 5402.2188
$ws.Range("K134").Value = 16206.6564
$ws.Range("M134").Value = -13671.6564

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H21").Value = 2000
$ws.Range("I21").Value = 2000
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 2000
$ws.Range("L21").Value = 0
$ws.Range("M21").Value = -1765
$ws.Range("N21").ClearContents()
$ws.Range("H31").Value = 3089823.8
$ws.Range("I31").Value = 3706788.5
$ws.Range("K31").Value = 3706788.5
$ws.Range("M31").Value = -3706493.5
$ws.Range("H34").Value = 3089823.8
$ws.Range("I34").Value = 3706788.5
$ws.Range("K34").Value = 3706788.5
$ws.Range("M34").Value = -3706586.5
$ws.Range("H132").Value = 1453.9231
$ws.Range("I132").Value = 1453.9231
$ws.Range("K132").Value = 4361.7693
$ws.Range("M132").Value = -1831.7693

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H52").Value = 3964.3333
$ws.Range("J52").Value = 3964.3333
$ws.Range("L52").Value = 11892.9999
$ws.Range("N52").Value = -12424.9999
$ws.Range("H68").Value = 7127.657
$ws.Range("J68").Value = 7970.484
$ws.Range("L68").Value = 23911.452
$ws.Range("N68").Value = -25533.452
$ws.Range("H71").Value = 7127.657
$ws.Range("J71").Value = 7970.484
$ws.Range("L71").Value = 71734.356
$ws.Range("N71").Value = -79846.356
$ws.Range("H117").Value = 8004.6665
$ws.Range("J117").Value = 16919.25
$ws.Range("L117").Value = 50757.75
$ws.Range("N117").Value = -57641.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 42000
$ws.Range("I15").Value = 42000
$ws.Range("K15").Value = 42000
$ws.Range("M15").Value = -41712
$ws.Range("H24").Value = 15627.25
$ws.Range("J24").Value = 19301.8
$ws.Range("L24").Value = 19301.8
$ws.Range("N24").Value = -19647.8
$ws.Range("H68").Value = 34971
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 34971
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 34971
$ws.Range("M68").ClearContents()
$ws.Range("N68").Value = -36593
$ws.Range("H70").Value = 68210.664
$ws.Range("I70").Value = 16868.6
$ws.Range("K70").Value = 16868.6
$ws.Range("M70").Value = -16598.6
$ws.Range("H71").Value = 34971
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 34971
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 104913
$ws.Range("M71").ClearContents()
$ws.Range("N71").Value = -113025
$ws.Range("H73").Value = 68210.664
$ws.Range("I73").Value = 16868.6
$ws.Range("K73").Value = 16868.6
$ws.Range("M73").Value = -15932.6
$ws.Range("H80").Value = 2612.5
$ws.Range("I80").Value = 2650
$ws.Range("K80").Value = 2650
$ws.Range("M80").Value = -1652
$ws.Range("H81").Value = 42000
$ws.Range("I81").Value = 42000
$ws.Range("K81").Value = 42000
$ws.Range("M81").Value = -41002
$ws.Range("H83").Value = 2612.5
$ws.Range("I83").Value = 2650
$ws.Range("K83").Value = 13250
$ws.Range("M83").Value = -8258
$ws.Range("H84").Value = 42000
$ws.Range("I84").Value = 42000
$ws.Range("K84").Value = 126000
$ws.Range("M84").Value = -121008
$ws.Range("H122").Value = 6410.36
$ws.Range("J122").Value = 1624.5
$ws.Range("L122").Value = 4873.5
$ws.Range("N122").Value = -9773.5
$ws.Range("H136").Value = 56045.613
$ws.Range("J136").Value = 56045.613
$ws.Range("L136").Value = 168136.839
$ws.Range("N136").Value = -173236.839

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H11").Value = 1999
$ws.Range("J11").Value = 1999
$ws.Range("L11").Value = 1999
$ws.Range("N11").Value = -2279
$ws.Range("H46").Value = 1986.4546
$ws.Range("I46").Value = 364.33334
$ws.Range("J46").Value = 5462.4287
$ws.Range("K46").Value = 364.33334
$ws.Range("L46").Value = 5462.4287
$ws.Range("M46").Value = -176.33334
$ws.Range("N46").Value = -5838.4287
$ws.Range("H55").Value = 614.13635
$ws.Range("I55").Value = 633.4286
$ws.Range("K55").Value = 633.4286
$ws.Range("M55").Value = -460.4286
$ws.Range("H68").Value = 3125
$ws.Range("J68").Value = 6500
$ws.Range("L68").Value = 6500
$ws.Range("N68").Value = -7998
$ws.Range("H71").Value = 3125
$ws.Range("J71").Value = 6500
$ws.Range("L71").Value = 32500
$ws.Range("N71").Value = -39988
$ws.Range("H93").Value = 1668.5
$ws.Range("I93").Value = 1446.7142
$ws.Range("J93").Value = 2056.625
$ws.Range("K93").Value = 1446.7142
$ws.Range("L93").Value = 2056.625
$ws.Range("M93").Value = -198.7141999999999
$ws.Range("N93").Value = -4552.625
$ws.Range("H136").Value = 6103629
$ws.Range("I136").Value = 7816566
$ws.Range("K136").Value = 23449698
$ws.Range("M136").Value = -23447148

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 805
$ws.Range("I100").Value = 930.55554
$ws.Range("J100").Value = 240
$ws.Range("K100").Value = 1861.11108
$ws.Range("L100").Value = 480
$ws.Range("M100").Value = -1320.11108
$ws.Range("N100").Value = -1562
$ws.Range("H113").Value = 900.575
$ws.Range("I113").Value = 945.5185
$ws.Range("K113").Value = 2836.5555
$ws.Range("M113").Value = -666.5554999999999
$ws.Range("H122").Value = 91908.766
$ws.Range("J122").Value = 226713
$ws.Range("L122").Value = 680139
$ws.Range("N122").Value = -685039
$ws.Range("H132").Value = 11114116
$ws.Range("I132").Value = 15153922
$ws.Range("J132").Value = 4648.5
$ws.Range("K132").Value = 45461766
$ws.Range("L132").Value = 13945.5
$ws.Range("M132").Value = -45459236
$ws.Range("N132").Value = -19005.5
$ws.Range("H136").Value = 6408382.5
$ws.Range("I136").Value = 1403409.4
$ws.Range("J136").Value = 28573264
$ws.Range("K136").Value = 4210228.199999999
$ws.Range("L136").Value = 85719792
$ws.Range("M136").Value = -4207678.199999999
$ws.Range("N136").Value = -85724892
